$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 1-7: existing songs get "ctrl clicked" again -> their "time ago"
# label (column G) advances, and a new (still-blank) time cell in column H
# is stamped with the same number format used by column I.
$timeAgo = @{
    1 = "5 hours ago"
    2 = "5 hours ago"
    3 = "5 hours ago"
    4 = "5 hours ago"
    5 = "5 hours ago"
    6 = "4 hours ago"
    7 = "5 hours ago"
}
foreach ($r in 1..7) {
    $ws.Range("G$r").Value = $timeAgo[$r]
    $ws.Range("H$r").NumberFormat = "h:mm"
}

# --- Row 8: "In My Blood" / Shawn Mendes -> "Yellow" / Coldplay / Parachutes
$ws.Range("D8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("B8").Value = "ADD"
$ws.Range("C8").Value = "Yellow"
$ws.Range("E8").Value = "Coldplay"
$ws.Range("F8").Value = "Parachutes"
$ws.Range("H8").NumberFormat = "h:mm"
$ws.Range("H8").Value = 0.18541666666666667

# --- Row 9: "Someone Like You" / Adele / 21 -> "head first" / Christian French
$ws.Range("D9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("B9").Value = "ADD"
$ws.Range("C9").Value = "head first"
$ws.Range("E9").Value = "Christian French"
$ws.Range("F9").Value = "head first"
$ws.Range("H9").NumberFormat = "h:mm"
$ws.Range("H9").Value = 0.1076388888888889

# --- Row 10: brand-new "ADD" row for "Flux" / Ellie Goulding / Brightest Blue
$ws.Range("B10").Value = "ADD"
$ws.Range("C10").Value = "Flux"
$ws.Range("E10").Value = "Ellie Goulding"
$ws.Range("F10").Value = "Brightest Blue"
$ws.Range("H10").NumberFormat = "h:mm"
$ws.Range("H10").Value = 0.16041666666666668


